$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.617.50'
$ws.Range("E2").Value = '  -2.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.897.20'
$ws.Range("E3").Value = '  -2.01%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.26'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.62'
$ws.Range("E6").Value = '  -3.49%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.896.07'
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.00'
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("E12").Value = '  -2.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").Value = '  -1.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.88'
$ws.Range("E14").Value = '  -2.90%  '
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.377.07'
$ws.Range("E16").Value = '  -2.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.609.96'
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.56'
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.894.42'
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '433.74'
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("E21").Value = '  -3.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.659'
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("E23").Value = '  -2.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.42'
$ws.Range("E24").Value = '  -1.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.97'
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.96'
$ws.Range("E27").Value = '  -11.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.01'
$ws.Range("E28").Value = '  -5.86%  '
$ws.Range("E29").Value = '  +5.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.99'
$ws.Range("E30").Value = '  -4.34%  '
$ws.Range("E31").Value = '  -4.11%  '
$ws.Range("E32").Value = '  -7.51%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.56'
$ws.Range("E35").Value = '  -3.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.959'
$ws.Range("E36").Value = '  -3.41%  '
$ws.Range("E37").Value = '  -3.88%  '
$ws.Range("E38").Value = '  -1.78%  '
$ws.Range("E39").Value = '  -5.48%  '
$ws.Range("E40").Value = '  -8.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.21'
$ws.Range("E41").Value = '  -3.57%  '
$ws.Range("E42").Value = '  -3.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.57'
$ws.Range("E43").Value = '  +1.36%  '
$ws.Range("E44").Value = '  -5.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.690.94'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.02'
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '348.13'
$ws.Range("E48").Value = '  -3.25%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  -1.59%  '
$ws.Range("E51").Value = '  -5.17%  '
